$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated TPM-derived numeric values for existing Sending-cluster rows (ECs, FAPs, Inflammatory-Mac)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5247883333333333
$ws.Range("H2").Value = 1.574365
$ws.Range("I2").Value = 0.1674845870648259
$ws.Range("J2").Value = 0.1674845870648259
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1011536666666667
$ws.Range("N2").Value = 0.303461
$ws.Range("O2").Value = 0.007629860605400263
$ws.Range("P2").Value = 0.008254451482408482
$ws.Range("Q2").Value = 0.05308426414055555
$ws.Range("R2").Value = 0.477758377265
$ws.Range("S2").Value = 0.001277884052857645
$ws.Range("T2").Value = 0.001382493397977825
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5247883333333333
$ws.Range("H3").Value = 1.574365
$ws.Range("I3").Value = 0.1674845870648259
$ws.Range("J3").Value = 0.1674845870648259
$ws.Range("O3").Value = 0.7653686681256785
$ws.Range("P3").Value = 0.8280227993585454
$ws.Range("Q3").Value = 5.325003252999444
$ws.Range("R3").Value = 47.925029276995
$ws.Range("S3").Value = 0.128187455333385
$ws.Range("T3").Value = 0.1386810566308271
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5247883333333333
$ws.Range("H4").Value = 1.574365
$ws.Range("I4").Value = 0.1674845870648259
$ws.Range("J4").Value = 0.1674845870648259
$ws.Range("M4").Value = 3.0094955
$ws.Range("N4").Value = 6.018991
$ws.Range("O4").Value = 0.2270014712689213
$ws.Range("P4").Value = 0.1637227491590462
$ws.Range("Q4").Value = 1.579348127619167
$ws.Range("R4").Value = 9.476088765715
$ws.Range("S4").Value = 0.03801924767858321
$ws.Range("T4").Value = 0.02742103703602092
$ws.Range("I5").Value = 0.7148668960482055
$ws.Range("J5").Value = 0.7148668960482057
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1011536666666667
$ws.Range("N5").Value = 0.303461
$ws.Range("O5").Value = 0.007629860605400263
$ws.Range("P5").Value = 0.008254451482408482
$ws.Range("Q5").Value = 0.2265771662945555
$ws.Range("R5").Value = 2.039194496651
$ws.Range("S5").Value = 0.005454334768262969
$ws.Range("T5").Value = 0.005900834109809862
$ws.Range("I6").Value = 0.7148668960482055
$ws.Range("J6").Value = 0.7148668960482057
$ws.Range("O6").Value = 0.7653686681256785
$ws.Range("P6").Value = 0.8280227993585454
$ws.Range("S6").Value = 0.5471367241155529
$ws.Range("T6").Value = 0.5919260884345896
$ws.Range("I7").Value = 0.7148668960482055
$ws.Range("J7").Value = 0.7148668960482057
$ws.Range("M7").Value = 3.0094955
$ws.Range("N7").Value = 6.018991
$ws.Range("O7").Value = 0.2270014712689213
$ws.Range("P7").Value = 0.1637227491590462
$ws.Range("Q7").Value = 6.741060258480166
$ws.Range("R7").Value = 40.446361550881
$ws.Range("S7").Value = 0.1622758371643896
$ws.Range("T7").Value = 0.1170399735038063
$ws.Range("G8").Value = 0.3686343333333333
$ws.Range("H8").Value = 1.105903
$ws.Range("I8").Value = 0.1176485168869685
$ws.Range("J8").Value = 0.1176485168869685
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1011536666666667
$ws.Range("N8").Value = 0.303461
$ws.Range("O8").Value = 0.007629860605400263
$ws.Range("P8").Value = 0.008254451482408482
$ws.Range("Q8").Value = 0.03728871447588888
$ws.Range("R8").Value = 0.335598430283
$ws.Range("S8").Value = 0.0008976417842796484
$ws.Range("T8").Value = 0.0009711239746207965
$ws.Range("G9").Value = 0.3686343333333333
$ws.Range("H9").Value = 1.105903
$ws.Range("I9").Value = 0.1176485168869685
$ws.Range("J9").Value = 0.1176485168869685
$ws.Range("O9").Value = 0.7653686681256785
$ws.Range("P9").Value = 0.8280227993585454
$ws.Range("Q9").Value = 3.740515746032111
$ws.Range("R9").Value = 33.664641714289
$ws.Range("S9").Value = 0.09004448867674046
$ws.Range("T9").Value = 0.09741565429312875
$ws.Range("G10").Value = 0.3686343333333333
$ws.Range("H10").Value = 1.105903
$ws.Range("I10").Value = 0.1176485168869685
$ws.Range("J10").Value = 0.1176485168869685
$ws.Range("M10").Value = 3.0094955
$ws.Range("N10").Value = 6.018991
$ws.Range("O10").Value = 0.2270014712689213
$ws.Range("P10").Value = 0.1637227491590462
$ws.Range("Q10").Value = 1.109403367312167
$ws.Range("R10").Value = 6.656420203873
$ws.Range("S10").Value = 0.02670638642594837
$ws.Range("T10").Value = 0.01926173861921895

# Remove rows for the "MuSCs" sending-cluster (no longer present after TPM recalculation)
$ws.Range("A11:T13").EntireRow.Delete()
